$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 327 number-format (A:L) as a style template for every new row,
# and copy the text-cell style (from L327) into whichever of M/N is used.

# ---- Row 328 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A328:L328").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N328").PasteSpecial(-4122)
$ws.Range("A328").Value = 45192.67438685185
$ws.Range("B328").Value = "gys5785@naver.com"
$ws.Range("C328").Value = "정치행정학과"
$ws.Range("D328").Value = 20232402
$ws.Range("E328").Value = "고형승"
$ws.Range("F328").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G328").Value = 0.9
$ws.Range("H328").Value = "4:6"
$ws.Range("I328").Value = "15분의 1"
$ws.Range("J328").Value = "130만호, 5백만명"
$ws.Range("K328").Value = "평안"
$ws.Range("L328").Value = "Black"
$ws.Range("N328").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# ---- Row 329 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A329:L329").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M329").PasteSpecial(-4122)
$ws.Range("A329").Value = 45192.68013449074
$ws.Range("B329").Value = "minjoo902@naver.com"
$ws.Range("C329").Value = "금융재무학과"
$ws.Range("D329").Value = 20192827
$ws.Range("E329").Value = "김민주"
$ws.Range("F329").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G329").Value = 0.1
$ws.Range("H329").Value = "6:4"
$ws.Range("I329").Value = "20분의 1"
$ws.Range("J329").Value = "20만호, 69만명"
$ws.Range("K329").Value = "충청"
$ws.Range("L329").Value = "Red"
$ws.Range("M329").Value = "모름/무응답"

# ---- Row 330 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A330:L330").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N330").PasteSpecial(-4122)
$ws.Range("A330").Value = 45192.6832569213
$ws.Range("B330").Value = "tkdgjs9768@naver.com"
$ws.Range("C330").Value = "경제학과"
$ws.Range("D330").Value = 20212837
$ws.Range("E330").Value = "임상헌"
$ws.Range("F330").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G330").Value = 0.1
$ws.Range("H330").Value = "6:4"
$ws.Range("I330").Value = "20분의 1"
$ws.Range("J330").Value = "20만호, 69만명"
$ws.Range("K330").Value = "충청"
$ws.Range("L330").Value = "Black"
$ws.Range("N330").Value = "찬성한다."

# ---- Row 331 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A331:L331").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N331").PasteSpecial(-4122)
$ws.Range("A331").Value = 45192.68467541666
$ws.Range("B331").Value = "seongmo0731@naver.com"
$ws.Range("C331").Value = "경영학과"
$ws.Range("D331").Value = 20192988
$ws.Range("E331").Value = "조성모"
$ws.Range("F331").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G331").Value = 0.1
$ws.Range("H331").Value = "6:4"
$ws.Range("I331").Value = "10분의 1"
$ws.Range("J331").Value = "20만호, 69만명"
$ws.Range("K331").Value = "충청"
$ws.Range("L331").Value = "Black"
$ws.Range("N331").Value = "모름/무응답"

# ---- Row 332 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A332:L332").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N332").PasteSpecial(-4122)
$ws.Range("A332").Value = 45192.689555590274
$ws.Range("B332").Value = "1004soeun@naver.com"
$ws.Range("C332").Value = "미디어스쿨"
$ws.Range("D332").Value = 20232514
$ws.Range("E332").Value = "김소은"
$ws.Range("F332").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G332").Value = 0.1
$ws.Range("H332").Value = "6:4"
$ws.Range("I332").Value = "10분의 1"
$ws.Range("J332").Value = "20만호, 69만명"
$ws.Range("K332").Value = "충청"
$ws.Range("L332").Value = "Black"
$ws.Range("N332").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# ---- Row 333 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A333:L333").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N333").PasteSpecial(-4122)
$ws.Range("A333").Value = 45192.691171493054
$ws.Range("B333").Value = "kgy5988@naver.com"
$ws.Range("C333").Value = "소프트웨어학부"
$ws.Range("D333").Value = 20203214
$ws.Range("E333").Value = "김진범"
$ws.Range("F333").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G333").Value = 0.1
$ws.Range("H333").Value = "6:4"
$ws.Range("I333").Value = "20분의 1"
$ws.Range("J333").Value = "20만호, 69만명"
$ws.Range("K333").Value = "충청"
$ws.Range("L333").Value = "Black"
$ws.Range("N333").Value = "찬성한다."

# ---- Row 334 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A334:L334").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M334").PasteSpecial(-4122)
$ws.Range("A334").Value = 45192.69671621528
$ws.Range("B334").Value = "sysy050300@naver.com"
$ws.Range("C334").Value = "심리학과"
$ws.Range("D334").Value = 20232101
$ws.Range("E334").Value = "고서연"
$ws.Range("F334").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G334").Value = 0.1
$ws.Range("H334").Value = "6:4"
$ws.Range("I334").Value = "20분의 1"
$ws.Range("J334").Value = "20만호, 69만명"
$ws.Range("K334").Value = "충청"
$ws.Range("L334").Value = "Red"
$ws.Range("M334").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# ---- Row 335 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A335:L335").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M335").PasteSpecial(-4122)
$ws.Range("A335").Value = 45192.69886010417
$ws.Range("B335").Value = "ggr1042@naver.com"
$ws.Range("C335").Value = "중국학과"
$ws.Range("D335").Value = 20221542
$ws.Range("E335").Value = "김경록"
$ws.Range("F335").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G335").Value = 0.5
$ws.Range("H335").Value = "3:7"
$ws.Range("I335").Value = "10분의 1"
$ws.Range("J335").Value = "20만호, 69만명"
$ws.Range("K335").Value = "충청"
$ws.Range("L335").Value = "Red"
$ws.Range("M335").Value = "반대한다."

# ---- Row 336 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A336:L336").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M336").PasteSpecial(-4122)
$ws.Range("A336").Value = 45192.702279074074
$ws.Range("B336").Value = "ntkrud0221@naver.com"
$ws.Range("C336").Value = "체육학과"
$ws.Range("D336").Value = 20234118
$ws.Range("E336").Value = "노태경"
$ws.Range("F336").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G336").Value = 0.1
$ws.Range("H336").Value = "6:4"
$ws.Range("I336").Value = "15분의 1"
$ws.Range("J336").Value = "20만호, 69만명"
$ws.Range("K336").Value = "경기"
$ws.Range("L336").Value = "Red"
$ws.Range("M336").Value = "모름/무응답"

# ---- Row 337 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A337:L337").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N337").PasteSpecial(-4122)
$ws.Range("A337").Value = 45192.70334378473
$ws.Range("B337").Value = "at79711@naver.com"
$ws.Range("C337").Value = "데이터사이언스"
$ws.Range("D337").Value = 20233220
$ws.Range("E337").Value = "박재영"
$ws.Range("F337").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G337").Value = 0.5
$ws.Range("H337").Value = "4:6"
$ws.Range("I337").Value = "10분의 1"
$ws.Range("J337").Value = "44만호, 153만명"
$ws.Range("K337").Value = "평안"
$ws.Range("L337").Value = "Black"
$ws.Range("N337").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# ---- Row 338 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A338:L338").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N338").PasteSpecial(-4122)
$ws.Range("A338").Value = 45192.71259708333
$ws.Range("B338").Value = "p51008085@gmail.com"
$ws.Range("C338").Value = "인공지능융합학부"
$ws.Range("D338").Value = 20236726
$ws.Range("E338").Value = "박준수"
$ws.Range("F338").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G338").Value = 0.1
$ws.Range("H338").Value = "6:4"
$ws.Range("I338").Value = "20분의 1"
$ws.Range("J338").Value = "20만호, 69만명"
$ws.Range("K338").Value = "충청"
$ws.Range("L338").Value = "Black"
$ws.Range("N338").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# ---- Row 339 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A339:L339").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M339").PasteSpecial(-4122)
$ws.Range("A339").Value = 45192.71808814815
$ws.Range("B339").Value = "jaejae7070@naver.com"
$ws.Range("C339").Value = "일본학과"
$ws.Range("D339").Value = 20221631
$ws.Range("E339").Value = "이재빈"
$ws.Range("F339").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G339").Value = 0.1
$ws.Range("H339").Value = "7:3"
$ws.Range("I339").Value = "10분의 1"
$ws.Range("J339").Value = "15만호,  32만명"
$ws.Range("K339").Value = "경기"
$ws.Range("L339").Value = "Red"
$ws.Range("M339").Value = "반대한다."

# ---- Row 340 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A340:L340").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N340").PasteSpecial(-4122)
$ws.Range("A340").Value = 45192.71840479167
$ws.Range("B340").Value = "leedongbin01@naver.com"
$ws.Range("C340").Value = "영어영문학과"
$ws.Range("D340").Value = 20231224
$ws.Range("E340").Value = "이동빈"
$ws.Range("F340").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G340").Value = 0.7
$ws.Range("H340").Value = "3:7"
$ws.Range("I340").Value = "10분의 1"
$ws.Range("J340").Value = "20만호, 69만명"
$ws.Range("K340").Value = "평안"
$ws.Range("L340").Value = "Black"
$ws.Range("N340").Value = "찬성한다."

# ---- Row 341 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A341:L341").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M341").PasteSpecial(-4122)
$ws.Range("A341").Value = 45192.72684841436
$ws.Range("B341").Value = "20217134@hallym.ac.kr"
$ws.Range("C341").Value = "체육학과"
$ws.Range("D341").Value = 20217134
$ws.Range("E341").Value = "장효경"
$ws.Range("F341").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G341").Value = 0.1
$ws.Range("H341").Value = "6:4"
$ws.Range("I341").Value = "10분의 1"
$ws.Range("J341").Value = "20만호, 69만명"
$ws.Range("K341").Value = "평안"
$ws.Range("L341").Value = "Red"
$ws.Range("M341").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# ---- Row 342 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A342:L342").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N342").PasteSpecial(-4122)
$ws.Range("A342").Value = 45192.728155578705
$ws.Range("B342").Value = "jinwoo3817@naver.com"
$ws.Range("C342").Value = "디지털미디어콘텐츠"
$ws.Range("D342").Value = 20222552
$ws.Range("E342").Value = "원진우"
$ws.Range("F342").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G342").Value = 0.1
$ws.Range("H342").Value = "3:7"
$ws.Range("I342").Value = "10분의 1"
$ws.Range("J342").Value = "20만호, 69만명"
$ws.Range("K342").Value = "충청"
$ws.Range("L342").Value = "Black"
$ws.Range("N342").Value = "모름/무응답"

# ---- Row 343 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A343:L343").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M343").PasteSpecial(-4122)
$ws.Range("A343").Value = 45192.72966489583
$ws.Range("B343").Value = "han7434370@naver.com"
$ws.Range("C343").Value = "체육학과"
$ws.Range("D343").Value = 20224152
$ws.Range("E343").Value = "한진우"
$ws.Range("F343").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G343").Value = 0.1
$ws.Range("H343").Value = "4:6"
$ws.Range("I343").Value = "10분의 1"
$ws.Range("J343").Value = "20만호, 69만명"
$ws.Range("K343").Value = "전라"
$ws.Range("L343").Value = "Red"
$ws.Range("M343").Value = "반대한다."

# ---- Row 344 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A344:L344").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N344").PasteSpecial(-4122)
$ws.Range("A344").Value = 45192.73511366898
$ws.Range("B344").Value = "qasw0529@naver.com"
$ws.Range("C344").Value = "미디어스쿨"
$ws.Range("D344").Value = 20232508
$ws.Range("E344").Value = "김민서"
$ws.Range("F344").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G344").Value = 0.1
$ws.Range("H344").Value = "3:7"
$ws.Range("I344").Value = "10분의 1"
$ws.Range("J344").Value = "20만호, 69만명"
$ws.Range("K344").Value = "전라"
$ws.Range("L344").Value = "Black"
$ws.Range("N344").Value = "찬성한다."

# ---- Row 345 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A345:L345").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N345").PasteSpecial(-4122)
$ws.Range("A345").Value = 45192.74084903936
$ws.Range("B345").Value = "jamesjm0612@gmail.com"
$ws.Range("C345").Value = "영어영문학과"
$ws.Range("D345").Value = 20231231
$ws.Range("E345").Value = "정재민"
$ws.Range("F345").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G345").Value = 0.3
$ws.Range("H345").Value = "6:4"
$ws.Range("I345").Value = "15분의 1"
$ws.Range("J345").Value = "20만호, 69만명"
$ws.Range("K345").Value = "평안"
$ws.Range("L345").Value = "Black"
$ws.Range("N345").Value = "모름/무응답"

# ---- Row 346 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A346:L346").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N346").PasteSpecial(-4122)
$ws.Range("A346").Value = 45192.74926520833
$ws.Range("B346").Value = "meldek98@gmail.com"
$ws.Range("C346").Value = "소프트웨어학부"
$ws.Range("D346").Value = 20235102
$ws.Range("E346").Value = "강비성"
$ws.Range("F346").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G346").Value = 0.1
$ws.Range("H346").Value = "6:4"
$ws.Range("I346").Value = "20분의 1"
$ws.Range("J346").Value = "20만호, 69만명"
$ws.Range("K346").Value = "충청"
$ws.Range("L346").Value = "Black"
$ws.Range("N346").Value = "찬성한다."

# ---- Row 347 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A347:L347").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N347").PasteSpecial(-4122)
$ws.Range("A347").Value = 45192.76095081019
$ws.Range("B347").Value = "alscoco100@gmail.com"
$ws.Range("C347").Value = "식품영양학과"
$ws.Range("D347").Value = 20223806
$ws.Range("E347").Value = "김민채"
$ws.Range("F347").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G347").Value = 0.7
$ws.Range("H347").Value = "5:5"
$ws.Range("I347").Value = "20분의 1"
$ws.Range("J347").Value = "15만호,  32만명"
$ws.Range("K347").Value = "충청"
$ws.Range("L347").Value = "Black"
$ws.Range("N347").Value = "모름/무응답"

# ---- Row 348 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A348:L348").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M348").PasteSpecial(-4122)
$ws.Range("A348").Value = 45192.767777824076
$ws.Range("B348").Value = "yeshin05@naver.com"
$ws.Range("C348").Value = "미래융합스쿨"
$ws.Range("D348").Value = 20236639
$ws.Range("E348").Value = "최예원"
$ws.Range("F348").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G348").Value = 0.1
$ws.Range("H348").Value = "6:4"
$ws.Range("I348").Value = "20분의 1"
$ws.Range("J348").Value = "20만호, 69만명"
$ws.Range("K348").Value = "충청"
$ws.Range("L348").Value = "Red"
$ws.Range("M348").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# ---- Row 349 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A349:L349").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M349").PasteSpecial(-4122)
$ws.Range("A349").Value = 45192.77363125
$ws.Range("B349").Value = "eugene3551@gmail.com"
$ws.Range("C349").Value = "소프트웨어학부"
$ws.Range("D349").Value = 20235214
$ws.Range("E349").Value = "유수영"
$ws.Range("F349").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G349").Value = 0.9
$ws.Range("H349").Value = "6:4"
$ws.Range("I349").Value = "30분의 1"
$ws.Range("J349").Value = "20만호, 69만명"
$ws.Range("K349").Value = "충청"
$ws.Range("L349").Value = "Red"
$ws.Range("M349").Value = "반대한다."

# ---- Row 350 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A350:L350").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("N350").PasteSpecial(-4122)
$ws.Range("A350").Value = 45192.78515171296
$ws.Range("B350").Value = "cozyandrelaxing2@gmail.com"
$ws.Range("C350").Value = "환경생명공학과"
$ws.Range("D350").Value = 20223725
$ws.Range("E350").Value = "이성민"
$ws.Range("F350").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G350").Value = 0.7
$ws.Range("H350").Value = "4:6"
$ws.Range("I350").Value = "10분의 1"
$ws.Range("J350").Value = "20만호, 69만명"
$ws.Range("K350").Value = "전라"
$ws.Range("L350").Value = "Black"
$ws.Range("N350").Value = "모름/무응답"

# ---- Row 351 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A351:L351").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M351").PasteSpecial(-4122)
$ws.Range("A351").Value = 45192.818970231485
$ws.Range("B351").Value = "hyunbin7379@gmail.com"
$ws.Range("C351").Value = "경영학부"
$ws.Range("D351").Value = 20233036
$ws.Range("E351").Value = "정현빈"
$ws.Range("F351").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G351").Value = 0.5
$ws.Range("H351").Value = "5:5"
$ws.Range("I351").Value = "15분의 1"
$ws.Range("J351").Value = "130만호, 5백만명"
$ws.Range("K351").Value = "경기"
$ws.Range("L351").Value = "Red"
$ws.Range("M351").Value = "모름/무응답"

# ---- Row 352 ----
$ws.Range("A327:L327").Copy()
$ws.Range("A352:L352").PasteSpecial(-4122)
$ws.Range("L327").Copy()
$ws.Range("M352").PasteSpecial(-4122)
$ws.Range("A352").Value = 45192.82088186343
$ws.Range("B352").Value = "jyn10131@naver.com"
$ws.Range("C352").Value = "식품영양학과"
$ws.Range("D352").Value = 20233847
$ws.Range("E352").Value = "정예나"
$ws.Range("F352").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G352").Value = 0.3
$ws.Range("H352").Value = "6:4"
$ws.Range("I352").Value = "15분의 1"
$ws.Range("J352").Value = "44만호, 153만명"
$ws.Range("K352").Value = "경상"
$ws.Range("L352").Value = "Red"
$ws.Range("M352").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

$excel.Application.CutCopyMode = $false
